# Apply the ALC/ARM/BSM/CRP/CUL/GSM/LTW price & profit updates scraped
# from the upstream price-refresh run. Values/rows identified by
# cross-referencing the canonical-OOXML diff against each sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H19").Value = 840.9167
$ws.Range("I19").Value = 759.8
$ws.Range("J19").Value = 898.8570999999999
$ws.Range("K19").Value = 759.8
$ws.Range("L19").Value = 898.8570999999999
$ws.Range("M19").Value = -584.8
$ws.Range("N19").Value = -1248.8571

$ws.Range("H39").Value = 704.3
$ws.Range("I39").Value = 858.2857
$ws.Range("J39").Value = 345
$ws.Range("K39").Value = 2574.8571
$ws.Range("L39").Value = 1035
$ws.Range("M39").Value = -2278.8571
$ws.Range("N39").Value = -1627

$ws.Range("H40").Value = 2038
$ws.Range("I40").Value = 1750
$ws.Range("J40").Value = 2326
$ws.Range("K40").Value = 1750
$ws.Range("L40").Value = 2326
$ws.Range("M40").Value = -1575
$ws.Range("N40").Value = -2676

$ws.Range("H43").Value = 860.4
$ws.Range("I43").Value = 950
$ws.Range("J43").Value = 800.6667
$ws.Range("K43").Value = 950
$ws.Range("L43").Value = 800.6667
$ws.Range("M43").Value = -881
$ws.Range("N43").Value = -938.6667

$ws.Range("H98").Value = 320379.7
$ws.Range("I98").Value = 399620.6
$ws.Range("J98").Value = 3416
$ws.Range("K98").Value = 399620.6
$ws.Range("L98").Value = 3416
$ws.Range("M98").Value = -398122.6
$ws.Range("N98").Value = -6412

$ws.Range("H116").Value = 19775068
$ws.Range("I116").Value = 46124164
$ws.Range("J116").Value = 13247.25
$ws.Range("K116").Value = 46124164
$ws.Range("L116").Value = 13247.25
$ws.Range("M116").Value = -46120722
$ws.Range("N116").Value = -20131.25

$ws.Range("H122").Value = 320379.7
$ws.Range("I122").Value = 399620.6
$ws.Range("J122").Value = 3416
$ws.Range("K122").Value = 1198861.8
$ws.Range("L122").Value = 10248
$ws.Range("M122").Value = -1196411.8
$ws.Range("N122").Value = -15148

$ws.Range("H131").Value = 6864.533
$ws.Range("I131").Value = 1497.5
$ws.Range("J131").Value = 12998.286
$ws.Range("K131").Value = 4492.5
$ws.Range("L131").Value = 38994.858
$ws.Range("M131").Value = 547.5
$ws.Range("N131").Value = -49074.858

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H61").Value = 2603.3
$ws.Range("I61").Value = 1814.6875
$ws.Range("K61").Value = 1814.6875
$ws.Range("M61").Value = -1602.6875

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0

$ws.Range("H132").Value = 2657.6223
$ws.Range("I132").Value = 2005.4839
$ws.Range("J132").Value = 4101.643
$ws.Range("K132").Value = 6016.4517
$ws.Range("L132").Value = 12304.929
$ws.Range("M132").Value = -3486.4517
$ws.Range("N132").Value = -17364.929

$ws.Range("H136").Value = 2603.3
$ws.Range("I136").Value = 1814.6875
$ws.Range("K136").Value = 5444.0625
$ws.Range("M136").Value = -2894.0625

$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H134").Value = 2069.2
$ws.Range("I134").Value = 1352.2572
$ws.Range("J134").Value = 4578.5
$ws.Range("K134").Value = 4056.7716
$ws.Range("L134").Value = 13735.5
$ws.Range("M134").Value = -1521.7716
$ws.Range("N134").Value = -18805.5

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H19").Value = 921.8182
$ws.Range("I19").Value = 921.8182
$ws.Range("K19").Value = 921.8182
$ws.Range("M19").Value = -751.8182

$ws.Range("H24").Value = 921.8182
$ws.Range("I24").Value = 921.8182
$ws.Range("K24").Value = 921.8182
$ws.Range("M24").Value = -751.8182

$ws.Range("H31").Value = 1849.0698
$ws.Range("I31").Value = 1088.3182
$ws.Range("J31").Value = 2646.0476
$ws.Range("K31").Value = 1088.3182
$ws.Range("L31").Value = 2646.0476
$ws.Range("M31").Value = -793.3181999999999
$ws.Range("N31").Value = -3236.0476

$ws.Range("H34").Value = 1849.0698
$ws.Range("I34").Value = 1088.3182
$ws.Range("J34").Value = 2646.0476
$ws.Range("K34").Value = 1088.3182
$ws.Range("L34").Value = 2646.0476
$ws.Range("M34").Value = -886.3181999999999
$ws.Range("N34").Value = -3050.0476

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 27716.455
$ws.Range("I4").Value = 268.7143
$ws.Range("J4").Value = 75750
$ws.Range("K4").Value = 806.1428999999999
$ws.Range("L4").Value = 227250
$ws.Range("M4").Value = -694.1428999999999
$ws.Range("N4").Value = -227474

$ws.Range("H5").Value = 1551.6842
$ws.Range("I5").Value = 1133.3334
$ws.Range("K5").Value = 3400.0002
$ws.Range("M5").Value = -3288.0002

$ws.Range("H129").Value = 2059
$ws.Range("I129").Value = 5245
$ws.Range("J129").Value = 1262.5
$ws.Range("K129").Value = 15735
$ws.Range("L129").Value = 3787.5
$ws.Range("M129").Value = -10735
$ws.Range("N129").Value = -13787.5

$ws.Range("H131").Value = 964.76
$ws.Range("I131").Value = 512.8570999999999
$ws.Range("J131").Value = 998.7742
$ws.Range("K131").Value = 1538.5713
$ws.Range("L131").Value = 2996.3226
$ws.Range("M131").Value = 3501.4287
$ws.Range("N131").Value = -13076.3226

$ws.Range("H135").Value = 1551.6842
$ws.Range("I135").Value = 1133.3334
$ws.Range("K135").Value = 10200.0006
$ws.Range("M135").Value = -7665.000599999999

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H9").Value = 838.6667
$ws.Range("I9").Value = 406.4
$ws.Range("K9").Value = 406.4
$ws.Range("M9").Value = -236.4

$ws.Range("H132").Value = 2952.0588
$ws.Range("I132").Value = 2340.3572
$ws.Range("J132").Value = 5806.6665
$ws.Range("K132").Value = 7021.071599999999
$ws.Range("L132").Value = 17419.9995
$ws.Range("M132").Value = -4491.071599999999
$ws.Range("N132").Value = -22479.9995

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 1017.1818
$ws.Range("I22").Value = 1021
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1021
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -726
$ws.Range("N22").Value = -1590

$ws.Range("H27").Value = 1017.1818
$ws.Range("I27").Value = 1021
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 1021
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -914
$ws.Range("N27").Value = -1214

$ws.Range("H55").Value = 216.16667
$ws.Range("I55").Value = 118.833336
$ws.Range("K55").Value = 118.833336
$ws.Range("M55").Value = 54.166664
